$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Usuarios")
$ws.Cells.Item(15, 1).Value = 0
$ws.Cells.Item(15, 2).Value = "Miguel"
$ws.Cells.Item(15, 3).Value = "Estudiante"
